# alteração da escala do sistema solar
#
# Adds a new "Raio do sol" row (D2:E2) above the existing table and rescales
# the virtual-solar-system numbers: F4 (scale factor) 500 -> 1000, E5 (real
# sun-to-planet distance) 384400 -> 57910000 (Mercury instead of the Moon),
# and the F5 formula now also folds in the sun's own radius (E2) alongside
# the planet radius (E5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 2: "Raio do sol" / 676000 -------------------------------------
$ws.Range("D2").Value = "Raio do sol"
$ws.Range("E2").Value = 676000
# Nudge the font so a dedicated cell-format record is minted for E2 (as in
# the authored workbook), while keeping the same visual (default) font.
$ws.Range("E2").Font.ThemeColor = 1

# --- Updated scale factor and real-world distance ---------------------------
$ws.Range("F4").Value = 1000
$ws.Range("E5").Value = 57910000

# --- Updated formula: include the sun's radius term -------------------------
$ws.Range("F5").Formula = "=((E5*F4)/E4) + ((E2*F4)/E4)"

# --- Column widths ------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 32.33203125
$ws.Columns.Item(5).ColumnWidth = 19
$ws.Columns.Item(6).ColumnWidth = 22

# --- Selection moves to F15 (as captured in the saved workbook state) -------
$ws.Range("F15").Select()
